# Update country data files
# Adds the MSME size-classification table (rows 18-22) and moves the
# "Source" attribution (SNC / long citation) down to rows 27-28 to make
# room for it, on the "Summary" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row for the new classification table (bold "title" look) ---
$ws.Range("B18").Value = "Number of employees"
$ws.Range("C18").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D18").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B18:D18").Font.Bold = $true

# --- Micro row ---
$ws.Range("A19").Value = "Micro"

# --- Small row ---
$ws.Range("A20").Value = "Small"
$ws.Range("B20").Value = "<50"
$ws.Range("D20").Value = "<100,000 UT"

# --- Medium row (overwrites the old A21 "SNC" cell; clear its old bold/title format) ---
$ws.Range("A21").ClearFormats()
$ws.Range("A21").Value = "Medium"
$ws.Range("B21").Value = "51-100 <br/><250 Industry, <br/><500 Trade, <br/><100 Service, <br/><50 Agriculture"
$ws.Range("D21").Value = "100,000 UT to 250,000 UT <br/><750,000 Industry, <br/><1,000,000 Trade, <br/><500,000 Serv., <br/><300,000 Agriculture"

# --- Large row (overwrites the old A22 "Servicio..." cell; clear its old italic/source format) ---
$ws.Range("A22").ClearFormats()
$ws.Range("A22").Value = "Large"
$ws.Range("B22").Value = ">100 <br/>>=250 Industry, <br/>>=500 Trade, <br/>>=100 Service, <br/>>=50 Agriculture"
$ws.Range("D22").Value = ">250,000 UT <br/>>=750,000 Industry, <br/>>=1,000,000 Trade, <br/>>=500,000 Serv.,<br/> >=300,000 Agriculture"

# --- Source attribution, moved down to rows 27-28 ---
$ws.Range("A27").Value = "SNC"
$ws.Range("A27").Font.Bold = $true

$ws.Range("A28").Value = "Servicio Nacional de Contrataciones (SNC), Foro de Estandares Internacionales, Papel del Contador Publico en la PYMEs Venezolanas. Available at http://fccpv.org/cont3/data/files/Foro-II-May2009-Presentacion-1.pdf"
$ws.Range("A28").Font.Italic = $true

Write-Host "Applied MSME classification table + moved source rows."
